# Auto-generated script applying value updates to match target diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 11999
$ws.Range("I7").Value = 11999
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 11999
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -11887
$ws.Range("N7").ClearContents()
$ws.Range("H14").Value = 11999
$ws.Range("I14").Value = 11999
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 11999
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -11808
$ws.Range("N14").ClearContents()
$ws.Range("H19").Value = 1494.4615
$ws.Range("I19").Value = 951.5
$ws.Range("J19").Value = 1735.7778
$ws.Range("K19").Value = 951.5
$ws.Range("L19").Value = 1735.7778
$ws.Range("M19").Value = -776.5
$ws.Range("N19").Value = -2085.7778
$ws.Range("H53").Value = 6678.7812
$ws.Range("I53").Value = 109.23077
$ws.Range("J53").Value = 11173.737
$ws.Range("K53").Value = 109.23077
$ws.Range("L53").Value = 11173.737
$ws.Range("M53").Value = 527.76923
$ws.Range("N53").Value = -12447.737
$ws.Range("H70").Value = 6129.478
$ws.Range("I70").Value = 5000
$ws.Range("J70").Value = 6154.5776
$ws.Range("K70").Value = 15000
$ws.Range("L70").Value = 18463.7328
$ws.Range("M70").Value = -14730
$ws.Range("N70").Value = -19003.7328
$ws.Range("H73").Value = 6129.478
$ws.Range("I73").Value = 5000
$ws.Range("J73").Value = 6154.5776
$ws.Range("K73").Value = 15000
$ws.Range("L73").Value = 18463.7328
$ws.Range("M73").Value = -14064
$ws.Range("N73").Value = -20335.7328
$ws.Range("H106").Value = 14126.25
$ws.Range("I106").Value = 15430
$ws.Range("K106").Value = 15430
$ws.Range("M106").Value = -14799
$ws.Range("H132").Value = 24393852
$ws.Range("I132").Value = 30307076
$ws.Range("J132").Value = 1800.125
$ws.Range("K132").Value = 90921228
$ws.Range("L132").Value = 5400.375
$ws.Range("M132").Value = -90918698
$ws.Range("N132").Value = -10460.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4881.18
$ws.Range("I32").Value = 3128.5938
$ws.Range("J32").Value = 7996.8887
$ws.Range("K32").Value = 3128.5938
$ws.Range("L32").Value = 7996.8887
$ws.Range("M32").Value = -2841.5938
$ws.Range("N32").Value = -8570.8887
$ws.Range("H61").Value = 16699.857
$ws.Range("I61").Value = 18650.666
$ws.Range("K61").Value = 18650.666
$ws.Range("M61").Value = -18438.666
$ws.Range("H102").Value = 3207939.5
$ws.Range("I102").Value = 3207939.5
$ws.Range("K102").Value = 3207939.5
$ws.Range("M102").Value = -3206317.5
$ws.Range("H122").Value = 698177.2
$ws.Range("I122").Value = 3428.6667
$ws.Range("J122").Value = 1392925.8
$ws.Range("K122").Value = 10286.0001
$ws.Range("L122").Value = 4178777.4
$ws.Range("M122").Value = -7836.000100000001
$ws.Range("N122").Value = -4183677.4
$ws.Range("H132").Value = 5533.579
$ws.Range("I132").Value = 5330.5625
$ws.Range("K132").Value = 15991.6875
$ws.Range("M132").Value = -13461.6875
$ws.Range("H136").Value = 16699.857
$ws.Range("I136").Value = 18650.666
$ws.Range("K136").Value = 55951.99800000001
$ws.Range("M136").Value = -53401.99800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4172147.5
$ws.Range("I86").Value = 6674659.5
$ws.Range("J86").Value = 1293.6666
$ws.Range("K86").Value = 6674659.5
$ws.Range("L86").Value = 1293.6666
$ws.Range("M86").Value = -6673536.5
$ws.Range("N86").Value = -3539.6666
$ws.Range("H89").Value = 4172147.5
$ws.Range("I89").Value = 6674659.5
$ws.Range("J89").Value = 1293.6666
$ws.Range("K89").Value = 33373297.5
$ws.Range("L89").Value = 6468.333000000001
$ws.Range("M89").Value = -33367681.5
$ws.Range("N89").Value = -17700.333
$ws.Range("H94").Value = 3500145.5
$ws.Range("I94").Value = 5050869.5
$ws.Range("K94").Value = 5050869.5
$ws.Range("M94").Value = -5050418.5
$ws.Range("H107").Value = 17858892
$ws.Range("I107").Value = 23811356
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 23811356
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = -23809436
$ws.Range("N107").Value = -5340
$ws.Range("H138").Value = 75320
$ws.Range("J138").Value = 75320
$ws.Range("L138").Value = 75320
$ws.Range("N138").Value = -85600
$ws.Range("H140").Value = 37999.668
$ws.Range("J140").Value = 37999.668
$ws.Range("L140").Value = 37999.668
$ws.Range("N140").Value = -48359.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 285.42856
$ws.Range("I7").Value = 193.66667
$ws.Range("K7").Value = 193.66667
$ws.Range("M7").Value = -80.66667000000001
$ws.Range("H16").Value = 2953.25
$ws.Range("I16").Value = 2407.4546
$ws.Range("K16").Value = 2407.4546
$ws.Range("M16").Value = -2120.4546
$ws.Range("H58").Value = 3390.55
$ws.Range("I58").Value = 3041.75
$ws.Range("J58").Value = 3913.75
$ws.Range("K58").Value = 3041.75
$ws.Range("L58").Value = 3913.75
$ws.Range("M58").Value = -2838.75
$ws.Range("N58").Value = -4319.75
$ws.Range("H105").Value = 1957
$ws.Range("I105").Value = 1889.125
$ws.Range("K105").Value = 1889.125
$ws.Range("M105").Value = -142.125
$ws.Range("H113").Value = 2953.25
$ws.Range("I113").Value = 2407.4546
$ws.Range("K113").Value = 2407.4546
$ws.Range("M113").Value = -237.4546
$ws.Range("H132").Value = 28310.236
$ws.Range("I132").Value = 34285.387
$ws.Range("J132").Value = 1848.8572
$ws.Range("K132").Value = 102856.161
$ws.Range("L132").Value = 5546.571599999999
$ws.Range("M132").Value = -100326.161
$ws.Range("N132").Value = -10606.5716
$ws.Range("H136").Value = 3390.55
$ws.Range("I136").Value = 3041.75
$ws.Range("J136").Value = 3913.75
$ws.Range("K136").Value = 9125.25
$ws.Range("L136").Value = 11741.25
$ws.Range("M136").Value = -6575.25
$ws.Range("N136").Value = -16841.25
$ws.Range("H138").Value = 70000
$ws.Range("J138").Value = 70000
$ws.Range("L138").Value = 70000
$ws.Range("N138").Value = -80280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1914.8334
$ws.Range("I14").Value = 1914.8334
$ws.Range("K14").Value = 5744.5002
$ws.Range("M14").Value = -5571.5002
$ws.Range("H92").Value = 1186.3334
$ws.Range("I92").Value = 930.25
$ws.Range("K92").Value = 2790.75
$ws.Range("M92").Value = -1542.75
$ws.Range("H97").Value = 450.2857
$ws.Range("I97").Value = 237.8
$ws.Range("K97").Value = 713.4000000000001
$ws.Range("M97").Value = -217.4000000000001
$ws.Range("H109").Value = 62504932
$ws.Range("I109").Value = 66671590
$ws.Range("K109").Value = 200014770
$ws.Range("M109").Value = -200013730
$ws.Range("H117").Value = 1281.0769
$ws.Range("I117").Value = 1599.5555
$ws.Range("J117").Value = 564.5
$ws.Range("K117").Value = 4798.666499999999
$ws.Range("L117").Value = 1693.5
$ws.Range("M117").Value = -1356.666499999999
$ws.Range("N117").Value = -8577.5
$ws.Range("H134").Value = 3208.6667
$ws.Range("I134").Value = 1394.6154
$ws.Range("K134").Value = 4183.8462
$ws.Range("M134").Value = 886.1538

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 10385.381
$ws.Range("I19").Value = 10121.647
$ws.Range("K19").Value = 10121.647
$ws.Range("M19").Value = -9833.647000000001
$ws.Range("H46").Value = 7149.8696
$ws.Range("J46").Value = 17777.777
$ws.Range("L46").Value = 17777.777
$ws.Range("N46").Value = -18089.777
$ws.Range("H57").Value = 17374.25
$ws.Range("J57").Value = 17856.285
$ws.Range("L57").Value = 17856.285
$ws.Range("N57").Value = -19496.285
$ws.Range("H133").Value = 109984
$ws.Range("J133").Value = 109984
$ws.Range("L133").Value = 109984
$ws.Range("N133").Value = -120104

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 37043704
$ws.Range("I61").Value = 111111110
$ws.Range("K61").Value = 111111110
$ws.Range("M61").Value = -111110908
$ws.Range("H113").Value = 37043704
$ws.Range("I113").Value = 111111110
$ws.Range("K113").Value = 111111110
$ws.Range("M113").Value = -111108940
$ws.Range("H122").Value = 6154.778
$ws.Range("J122").Value = 8824.571
$ws.Range("L122").Value = 26473.713
$ws.Range("N122").Value = -31373.713
$ws.Range("H132").Value = 7947.9697
$ws.Range("I132").Value = 7803.75
$ws.Range("K132").Value = 23411.25
$ws.Range("M132").Value = -20881.25
$ws.Range("H133").Value = 498998
$ws.Range("J133").Value = 498998
$ws.Range("L133").Value = 498998
$ws.Range("N133").Value = -504058

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 28703
$ws.Range("J94").Value = 28703
$ws.Range("L94").Value = 28703
$ws.Range("N94").Value = -30505
$ws.Range("H95").Value = 34810.75
$ws.Range("J95").Value = 34810.75
$ws.Range("L95").Value = 34810.75
$ws.Range("N95").Value = -40302.75
$ws.Range("H101").Value = 16166.333
$ws.Range("J101").Value = 16166.333
$ws.Range("L101").Value = 16166.333
$ws.Range("N101").Value = -22656.333
$ws.Range("H113").Value = 654.2727
$ws.Range("I113").Value = 614.625
$ws.Range("K113").Value = 1843.875
$ws.Range("M113").Value = 326.125
$ws.Range("H122").Value = 2433.1936
$ws.Range("J122").Value = 2470.1428
$ws.Range("L122").Value = 7410.428400000001
$ws.Range("N122").Value = -12310.4284
$ws.Range("H132").Value = 32617850
$ws.Range("I132").Value = 50006908
$ws.Range("K132").Value = 150020724
$ws.Range("M132").Value = -150018194
